# Weekly update: insert the latest week's price record for
# "Pepino ensalada" (Vega Central Mapocho de Santiago) at row 320,
# pushing all subsequent historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 320 (existing rows 320:392 shift to 321:393)
$ws.Rows("320:320").Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A320").Value = 9
$ws.Range("B320").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C320").Value = "Metropolitana"
$ws.Range("D320").Value = 44964
$ws.Range("E320").Value = 13
$ws.Range("F320").Value = 100112043
$ws.Range("G320").Value = "Pepino ensalada"
$ws.Range("H320").Value = "Sin especificar"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 70
$ws.Range("K320").Value = 9000
$ws.Range("L320").Value = 10000
$ws.Range("M320").Value = 9500
$ws.Range("N320").Value = "`$/caja 60 unidades"
$ws.Range("O320").Value = "Región Metropolitana"
$ws.Range("P320").Value = 158
$ws.Range("Q320").Value = 60
$ws.Range("R320").Value = "Hortaliza"
